# uso eficaz de pipes entre funciones y GUI
#
# Adds newly scraped catalogue rows to both sheets and normalizes a couple of
# cells that had been entered as text numbers into real numbers.

$wb = $excel.ActiveWorkbook

function Set-NumericText {
    <#
        Writes a genuine numeric value into a cell whose column style applies
        a Text ("@") number format. Assigning .Value directly while the "@"
        format is active makes Excel store the value as text, so we flip the
        cell to the default "Normal" style, write the number, then restore
        the Text display format (reusing the existing style record instead
        of minting a new one).
    #>
    param($cell, $value)

    $cell.Style = "Normal"
    $cell.Value = $value
    $cell.NumberFormat = "@"
}

# ---------------------------------------------------------------------------
# Sheet 1: "peliculas o documentales"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("peliculas o documentales")

# Row 10 - Sumergidos (Arriendo / Acción / HD / 2019 / $3.490)
$ws1.Range("A10").Value = "Sumergidos"
$ws1.Range("C10").Value = "Arriendo"
$ws1.Range("D10").Value = "Acción"
$ws1.Range("E10").Value = "HD"
Set-NumericText $ws1.Range("F10") 2019
$ws1.Range("G10").Value = "$3.490"

# Row 11 - Planeta Hostil (WILD / Premium / Documental / HD / 2019)
$ws1.Range("A11").Value = "Planeta Hostil"
$ws1.Range("B11").Value = "WILD"
$ws1.Range("C11").Value = "Premium"
$ws1.Range("D11").Value = "Documental"
$ws1.Range("E11").Value = "HD"
Set-NumericText $ws1.Range("F11") 2019

# Row 12 - This Is Us (FOX PREMIUM / Premium / Drama / HD / 2016)
$ws1.Range("A12").Value = "This Is Us"
$ws1.Range("B12").Value = "FOX PREMIUM"
$ws1.Range("C12").Value = "Premium"
$ws1.Range("D12").Value = "Drama"
$ws1.Range("E12").Value = "HD"
Set-NumericText $ws1.Range("F12") 2016

# Row 13 - Guerrilla del Oro (NATIONAL GEOGRAPHIC / Premium / Documental / HD / "2019" as text)
$ws1.Range("A13").Value = "Guerrilla del Oro"
$ws1.Range("B13").Value = "NATIONAL GEOGRAPHIC"
$ws1.Range("C13").Value = "Premium"
$ws1.Range("D13").Value = "Documental"
$ws1.Range("E13").Value = "HD"
$ws1.Range("F13").Value = "2019"

# ---------------------------------------------------------------------------
# Sheet 2: "series"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("series")

# Row 5 had its AGNO/EPISODIOS entered as text ("2020" / "1"); fix to real numbers.
Set-NumericText $ws2.Range("F5") 2020
Set-NumericText $ws2.Range("H5") 1

# Row 6 - duplicate of My Brilliant Friend / Gratis / Drama / HD / 2020 / Temporada 02 / 1
$ws2.Range("A6").Value = "My Brilliant Friend"
$ws2.Range("C6").Value = "Gratis"
$ws2.Range("D6").Value = "Drama"
$ws2.Range("E6").Value = "HD"
Set-NumericText $ws2.Range("F6") 2020
$ws2.Range("G6").Value = "Temporada 02"
Set-NumericText $ws2.Range("H6") 1

Write-Host "edit complete"
